$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '58.976.55'
$ws.Range("E2").Value = '  +1.41%  '
$ws.Range("D3").Value = '2.586.55'
$ws.Range("E3").Value = '  -0.26%  '
$ws.Range("D4").Value = '''1.00'
$ws.Range("E4").Value = '  -0.04%  '
$ws.Range("D5").Value = '''526.00'
$ws.Range("E5").Value = '  +0.74%  '
$ws.Range("D6").Value = '''138.93'
$ws.Range("E6").Value = '  -3.07%  '
$ws.Range("E7").Value = '  +0.07%  '
$ws.Range("D8").Value = '''0.563'
$ws.Range("E8").Value = '  -0.92%  '
$ws.Range("D9").Value = '2.596.75'
$ws.Range("E9").Value = '  -0.60%  '
$ws.Range("D10").Value = '''6.42'
$ws.Range("E10").Value = '  -0.63%  '
$ws.Range("E11").Value = '  -0.22%  '
$ws.Range("E12").Value = '  -3.21%  '
$ws.Range("E13").Value = '  +2.87%  '
$ws.Range("D14").Value = '3.047.43'
$ws.Range("D15").Value = '58.920.54'
$ws.Range("E15").Value = '  +1.38%  '
$ws.Range("D16").Value = '''20.49'
$ws.Range("E16").Value = '  +0.65%  '
$ws.Range("D17").Value = '2.582.03'
$ws.Range("E17").Value = '  -1.81%  '
$ws.Range("E18").Value = '  -0.96%  '
$ws.Range("D19").Value = '''343.61'
$ws.Range("E19").Value = '  +1.24%  '
$ws.Range("E20").Value = '  -0.84%  '
$ws.Range("D21").Value = '''10.05'
$ws.Range("E21").Value = '  -1.78%  '
$ws.Range("D22").Value = '''6.42'
$ws.Range("E22").Value = '  -0.09%  '
$ws.Range("D23").Value = '''0.998'
$ws.Range("E23").Value = '  +0.07%  '
$ws.Range("D24").Value = '''66.40'
$ws.Range("E25").Value = '  -0.16%  '
$ws.Range("E26").Value = '  +0.09%  '
$ws.Range("E27").Value = '  +0.03%  '
$ws.Range("D28").Value = '''7.05'
$ws.Range("E28").Value = '  +0.22%  '
$ws.Range("D30").Value = '0.0₃0720'
$ws.Range("E30").Value = '  -3.41%  '
$ws.Range("E31").Value = '  +1.34%  '
$ws.Range("D32").Value = '''5.89'
$ws.Range("E32").Value = '  -3.90%  '
$ws.Range("D33").Value = '''18.69'
$ws.Range("E33").Value = '  -0.45%  '
$ws.Range("D34").Value = '''149.54'
$ws.Range("E34").Value = '  -0.09%  '
$ws.Range("D35").Value = '''3.96'
$ws.Range("E35").Value = '  -1.11%  '
$ws.Range("E36").Value = '  -1.26%  '
$ws.Range("D37").Value = '''36.78'
$ws.Range("E37").Value = '  +2.15%  '
$ws.Range("E38").Value = '  +1.54%  '
$ws.Range("D39").Value = '''0.824'
$ws.Range("E39").Value = '  -4.46%  '
$ws.Range("D40").Value = '''0.807'
$ws.Range("E40").Value = '  -6.69%  '
$ws.Range("E41").Value = '  -0.55%  '
$ws.Range("D42").Value = '''0.997'
$ws.Range("D43").Value = '''0.601'
$ws.Range("E43").Value = '  -1.08%  '
$ws.Range("D44").Value = '''269.42'
$ws.Range("E44").Value = '  -0.57%  '
$ws.Range("E45").Value = '  +0.58%  '
$ws.Range("D46").Value = '''0.0953'
$ws.Range("E46").Value = '  -0.41%  '
$ws.Range("E47").Value = '  -1.62%  '
$ws.Range("D48").Value = '''18.36'
$ws.Range("E48").Value = '  -2.01%  '
$ws.Range("B49").Value = 'Maker'
$ws.Range("C49").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D49").Value = '1.963.00'
$ws.Range("E49").Value = '  -0.20%  '
$ws.Range("B50").Value = 'VeChain'
$ws.Range("C50").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D50").Value = '''0.0222'
$ws.Range("E50").Value = '  +0.05%  '
$ws.Range("B51").Value = 'InjectiveProtocol'
$ws.Range("C51").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D51").Value = '''18.15'
$ws.Range("E51").Value = '  -3.12%  '
